# Natmi following Dr Hou advice
#
# Rewrites the LR-pairs_lrc2p/Cthrc1-Ror2 sheet: the "Ligand-expressing
# cells" / "Receptor-expressing cells" counts (and every downstream
# statistic derived from them) were recomputed, and three new rows were
# added for the "sCs" -> {ECs,FAPs,sCs} target-cluster combinations that
# didn't exist before.
#
# Columns: A=Sending cluster, B=Ligand symbol, C=Receptor symbol,
#          D=Target cluster, E..T = the 16 numeric statistic columns.

$data = @(
    @("ECs",  "Cthrc1", "Ror2", "ECs",  3, 1, 0.6191033333333333, 1.85731,    0.007929637811697733, 0.007962490129789305, 1, 0.3333333333333333, 0.047572,           0.142716,   0.006780712238476238, 0.00730352994735907, 0.02945198377333334, 0.26506785396,      0.00005376859215646276, 0.0000581542851184672),
    @("ECs",  "Cthrc1", "Ror2", "FAPs", 3, 1, 0.6191033333333333, 1.85731,    0.007929637811697733, 0.007962490129789305, 3, 1,                  5.461551666666666, 16.384655,  0.7784665397132128,   0.8384891565742209,  3.381264842005555,   30.43138357805,     0.006172957708451387,   0.006676461633157593),
    @("ECs",  "Cthrc1", "Ror2", "sCs",  3, 1, 0.6191033333333333, 1.85731,    0.007929637811697733, 0.007962490129789305, 2, 1,                  1.5066585,         3.013317,   0.214752748048311,    0.15420731347842,    0.932777299545,      5.59666379727,      0.001702911511089883,   0.001227874211513245),
    @("FAPs", "Cthrc1", "Ror2", "ECs",  3, 1, 76.48912033333333,  229.467361, 0.9796927073757713,   0.9837515530909214,   1, 0.3333333333333333, 0.047572,           0.142716,   0.006780712238476238, 0.00730352994735907, 3.638740432497333,   32.748663892476,    0.006643014330848813,   0.007184858928760541),
    @("FAPs", "Cthrc1", "Ror2", "FAPs", 3, 1, 76.48912033333333,  229.467361, 0.9796927073757713,   0.9837515530909214,   3, 1,                  5.461551666666666, 16.384655,  0.7784665397132128,   0.8384891565742209,  417.7492826383838,   3759.743543745454,  0.7626579918930858,     0.8248650100297865),
    @("FAPs", "Cthrc1", "Ror2", "sCs",  3, 1, 76.48912033333333,  229.467361, 0.9796927073757713,   0.9837515530909214,   2, 1,                  1.5066585,         3.013317,   0.214752748048311,    0.15420731347842,    115.2429833077395,   691.4578998464369,  0.2103917011518367,     0.1517016841323743),
    @("sCs",  "Cthrc1", "Ror2", "ECs",  2, 1, 0.9663805,          1.932761,   0.01237765481253107,  0.008285956779289245, 1, 0.3333333333333333, 0.047572,           0.142716,   0.006780712238476238, 0.00730352994735907, 0.045972653146,      0.275835918876,     0.00008392931547096374, 0.00006051673348006191),
    @("sCs",  "Cthrc1", "Ror2", "FAPs", 2, 1, 0.9663805,          1.932761,   0.01237765481253107,  0.008285956779289245, 3, 1,                  5.461551666666666, 16.384655,  0.7784665397132128,   0.8384891565742209,  5.277937030409166,   31.66762218245499,  0.009635590111675658,   0.006947684911276686),
    @("sCs",  "Cthrc1", "Ror2", "sCs",  2, 1, 0.9663805,          1.932761,   0.01237765481253107,  0.008285956779289245, 2, 1,                  1.5066585,         3.013317,   0.214752748048311,    0.15420731347842,    1.45600539455925,    5.824021578237,     0.002658135385384449,   0.001277755134532496)
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Count; $j++) {
        $ws.Cells.Item($rowNum, $j + 1).Value = $rowVals[$j]
    }
}
